# Apply the "new review period" update: append a new date column
# (2022-04-05, serial 44656) with scores to the "3 Months", "12 Months"
# and "Summary" sheets of the Chicago strategy-meeting scores workbook.

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# "3 Months" sheet -> new column AB (28)
# ----------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("3 Months")
$col1 = 28

# Header date cell: value + matching date number-format (copy from the
# previous header cell so it reuses the same style as the rest of row 1).
$ws1.Cells.Item(1, $col1).Value = 44656
$ws1.Range("AA1").Copy()
$ws1.Cells.Item(1, $col1).PasteSpecial(-4122)

$scores1 = @{
    2  = -1
    3  = -2
    4  = -1
    5  = -2
    6  = -1
    7  = -1
    8  = -1
    9  = -2
    10 = -1
    11 = -1
    12 = -1
    13 = -1
    14 = -1
    15 = -1
    16 = -1
    17 = -1
    18 = -1
    19 = -1
    21 = -1
    22 = 0
    23 = 0
    24 = -2
    25 = -1
    26 = -2
}
foreach ($row in $scores1.Keys) {
    $ws1.Cells.Item($row, $col1).Value = $scores1[$row]
}

# ----------------------------------------------------------------------
# "12 Months" sheet -> new column AB (28)
# ----------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("12 Months")
$col2 = 28

$ws2.Cells.Item(1, $col2).Value = 44656
$ws2.Range("AA1").Copy()
$ws2.Cells.Item(1, $col2).PasteSpecial(-4122)

$scores2 = @{
    2  = -2
    3  = -2
    4  = -2
    5  = -1
    6  = -2
    7  = -1
    8  = -1
    9  = -2
    10 = -2
    11 = -2
    12 = -2
    13 = -1
    14 = 1
    15 = -2
    16 = -2
    17 = -1
    18 = -1
    19 = -1
    21 = 0
    22 = -2
    23 = 0
    24 = -1
    25 = -1
    26 = -2
}
foreach ($row in $scores2.Keys) {
    $ws2.Cells.Item($row, $col2).Value = $scores2[$row]
}

# ----------------------------------------------------------------------
# "Summary" sheet -> new column AY (51)
# ----------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Summary")
$col3 = 51

$ws3.Cells.Item(1, $col3).Value = 44656
$ws3.Range("AX1").Copy()
$ws3.Cells.Item(1, $col3).PasteSpecial(-4122)

$scores3 = @{
    2  = -1
    3  = -2
    4  = -1
    5  = -2
    6  = -2
    7  = -1
    8  = -1
    9  = -1
    10 = -1
    11 = 0
    12 = 0
    13 = -1
}
foreach ($row in $scores3.Keys) {
    $ws3.Cells.Item($row, $col3).Value = $scores3[$row]
}

# ----------------------------------------------------------------------
# Restore view state as closely as the host allows: active sheet +
# selected cells on each sheet (matches the saved selections in the diff).
# ----------------------------------------------------------------------
$ws1.Range("AD13").Select() | Out-Null
$ws2.Activate()
$ws2.Range("AB16").Select() | Out-Null
$ws3.Activate()
$ws3.Range("AY14").Select() | Out-Null
